$wb = $excel.ActiveWorkbook

# --- OFF sheet: update row 2 (Short Att, Short Comp, Deep Att, Deep Comp) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 209
$wsOff.Range("C2").Value = 139
$wsOff.Range("D2").Value = 63
$wsOff.Range("E2").Value = 23

# --- DEF sheet: update row 2 (Short Att, Short Comp, Deep Att, Deep Comp, Deep Int) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 158
$wsDef.Range("C2").Value = 111
$wsDef.Range("D2").Value = 35
$wsDef.Range("E2").Value = 20
$wsDef.Range("F2").Value = 3
